# Trading update: 2026-02-17 12:30:06
# Append a new OPEN MarketMaking trade (trade #21) as row 22 to both the
# "All Trades" log and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(22, 1).Value = 21                 # A22 Trade #

    # Date/Time are stored as plain text in this log (not Excel dates),
    # so force text formatting before assigning, then drop the format
    # again so no stray number-format style is left behind on the cell.
    $ws.Cells.Item(22, 2).NumberFormat = "@"
    $ws.Cells.Item(22, 2).Value = "2026-02-17"        # B22 Date
    $ws.Cells.Item(22, 2).ClearFormats()

    $ws.Cells.Item(22, 3).NumberFormat = "@"
    $ws.Cells.Item(22, 3).Value = "12:29:28"          # C22 Time
    $ws.Cells.Item(22, 3).ClearFormats()

    $ws.Cells.Item(22, 4).Value = "MarketMaking"      # D22 Strategy
    $ws.Cells.Item(22, 5).Value = "UP"                # E22 Side
    $ws.Cells.Item(22, 6).Value = 0.01                # F22 Entry Price
    # G22 Exit Price left blank - trade is still OPEN, no exit yet.
    $ws.Cells.Item(22, 8).Value = "OPEN"              # H22 Status
    $ws.Cells.Item(22, 9).Value = 0                   # I22 P&L %
    $ws.Cells.Item(22, 10).Value = 0                  # J22 P&L $
    $ws.Cells.Item(22, 11).Value = 99.99354434589566  # K22 Capital After
    $ws.Cells.Item(22, 12).Value = 0                  # L22 Entry Slippage (bps)
    $ws.Cells.Item(22, 13).Value = 0                  # M22 Exit Slippage (bps)
    $ws.Cells.Item(22, 14).Value = 0.6                # N22 Confidence
    $ws.Cells.Item(22, 15).Value = "Normal spread capture: 19600 bps" # O22 Entry Reason
    # P22 Exit Reason left blank - trade is still OPEN, no exit yet.
    $ws.Cells.Item(22, 17).Value = 0                  # Q22 Duration (min)
}
